$d = $word.ActiveDocument

# The document has a single section whose first-page / primary
# headers & footers carry the Pearson / BTec logo pictures as
# inline pictures. The logo picture names recorded in wp:docPr / pic:cNvPr
# need to be renamed:
#   - both Pearson logo footers: image2.png -> image1.png
#   - the BTec logo header:      image1.jpg -> image2.jpg

$sec = $d.Sections.Item(1)

# --- Footers: rename every Pearson logo inline picture ---
for ($f = 1; $f -le 3; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

# --- Headers: rename the BTec logo inline picture ---
for ($h = 1; $h -le 3; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($i = 1; $i -le $shapes.Count; $i++) {
            $shp = $shapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
